{"js": "// Fix a typo: \"\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" -> \"\u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" in the paragraph discussing\n// \"\u0414\u043b\u044f \u043f\u043e\u043d\u0438\u043c\u0430\u043d\u0438\u044f \u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e \u0442\u043e\u0433\u043e \u0444\u0430\u043a\u0442\u0430...\" (missing leading \"\u0434\").\n//\n// \"\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" also occurs as a substring inside \"\u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\u043c\" elsewhere in\n// the document, so we anchor on the unique surrounding phrase\n// \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" and replace it with \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\".\n\nconst results = context.document.body.search(\"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target phrase \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" not found.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix a typo: \"\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" -> \"\u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" in the paragraph discussing\n# \"\u0414\u043b\u044f \u043f\u043e\u043d\u0438\u043c\u0430\u043d\u0438\u044f \u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e \u0442\u043e\u0433\u043e \u0444\u0430\u043a\u0442\u0430...\" (missing leading \"\u0434\").\n#\n# \"\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" also occurs as a substring inside \"\u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\u043c\" elsewhere in\n# the document, so we anchor the Find on the unique surrounding phrase\n# \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\" and replace it with \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\".\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\",   # FindText\n    $true,                 # MatchCase\n    $false,                # MatchWholeWord\n    $false,                # MatchWildcards\n    $false,                # MatchSoundsLike\n    $false,                # MatchAllWordForms\n    $true,                 # Forward\n    $wdFindContinue,       # Wrap\n    $false,                # Format\n    \"\u043a\u0430\u0440\u0442\u0438\u043d\u044b \u0434\u043e\u0441\u0442\u0430\u0442\u043e\u0447\u043d\u043e\",  # ReplaceWith\n    $wdReplaceAll           # Replace\n)\n"}
